# Update the "dSF" column (F) with freshly recalculated delta-stock-final
# values for each row, reflecting the repulled data / mean calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new F value
$updates = @{
    2  = 1
    3  = 3
    4  = -1
    6  = 7
    7  = 9
    8  = -4
    9  = 4
    10 = 2
    11 = -1
    13 = 1
    14 = -3
    17 = -1
    18 = -1
    19 = 2
    22 = 1
    23 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
